$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos value cells changed to the docente name ---
$ws.Range("B10").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C10").Value = "4780627 - Ana Lucia Gabas Ferreira"

# --- Row 13: now carries the "Programa resumido:" label + "Semestral" value ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# --- Row 14: now "Short syllabus:" + the Unit operations paragraph ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Unit operations and processes: fluid rheology, sizing of pipes and fittings, pumping, stirring and mixing, characterization of particles and particle bed, sedimentation, filtration, processes with membrane. Unit operations of thermal exchange: heat exchangers and evaporators."
$ws.Range("C14").Value = "Unit operations and processes: fluid rheology, sizing of pipes and fittings, pumping, stirring and mixing, characterization of particles and particle bed, sedimentation, filtration, processes with membrane. Unit operations of thermal exchange: heat exchangers and evaporators."

# --- Row 15: now "Programa:" + the activation date ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"
$ws.Rows(15).RowHeight = 120

# --- Row 16: now "Syllabus:" + the English syllabus paragraph ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "- fluid rheology,- Sizing of pipes,- Accessories and pumping for industrial fluids,- Stirring and mixing,- Particle characterization and particle bed,- Sedimentation,- Filtration,- Processes with membranes.- Unit heat exchange operations: heat exchangers and evaporators."
$ws.Range("C16").Value = "- fluid rheology,- Sizing of pipes,- Accessories and pumping for industrial fluids,- Stirring and mixing,- Particle characterization and particle bed,- Sedimentation,- Filtration,- Processes with membranes.- Unit heat exchange operations: heat exchangers and evaporators."

# --- Row 17: now only "Avaliação:" label, value cells cleared, default height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows(17).RowHeight = $wb.ActiveSheet.StandardHeight

# --- Row 18: now "Método:" + the docente name ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C18").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Rows(18).RowHeight = 60

# --- Row 19: label becomes "Critério:" (value unchanged) ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: label becomes "Norma de recuperação:" (value unchanged) ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: label becomes "Bibliografia:" (value unchanged), height grows to 120 ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

# --- Row 22: now only "Requisitos:" label, value cells cleared, default height ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows(22).RowHeight = $wb.ActiveSheet.StandardHeight

# --- Row 23: label cleared, now carries the requisitos value text ---
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)`n"
$ws.Rows(23).RowHeight = 30

# --- Row 24 no longer exists in the sheet; remove it entirely ---
$ws.Rows(24).Delete()
